# Updated symbol list (cryptos) - refresh Price column and two Volume(1h) labels.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (column D) updates. A leading apostrophe is used so the
# numeric-looking values are stored as text, matching the workbook's
# existing inline-string convention for this column.
$ws.Range("D2").Value = "'248.57"
$ws.Range("D3").Value = "'21.81"
$ws.Range("D4").Value = "'5.354"
$ws.Range("D5").Value = "'0.05613"
$ws.Range("D6").Value = "'3.407"
$ws.Range("D7").Value = "'6.384"
$ws.Range("D8").Value = "'0.8144"
$ws.Range("D9").Value = "'0.9603"
$ws.Range("D10").Value = "'0.1416"
$ws.Range("D11").Value = "'0.07681"
$ws.Range("D13").Value = "'0.03054"
$ws.Range("D14").Value = "'0.09306"
$ws.Range("D15").Value = "'3.565"
$ws.Range("D16").Value = "'0.001609"
$ws.Range("D17").Value = "'0.04712"
$ws.Range("D18").Value = "'0.0005768"

# Volume(1h) label update for row 18 (One/ONE no longer flagged "Worst in 24h").
$ws.Range("E18").Value = "17OneONE"

$ws.Range("D19").Value = "'0.006458"
$ws.Range("D20").Value = "'0.005077"
$ws.Range("D21").Value = "'0.001032"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("D23").Value = "'3.751"
$ws.Range("D24").Value = "'2.143"
$ws.Range("D25").Value = "'0.3256"
$ws.Range("D28").Value = "'0.0003098"
$ws.Range("D40").Value = "'0.03952"
$ws.Range("D41").Value = "'0.006967"
$ws.Range("D42").Value = "'0.1062"
$ws.Range("D43").Value = "'0.003030"
$ws.Range("D44").Value = "'0.008622"
$ws.Range("D45").Value = "'0.00005815"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D47").Value = "'0.0005498"

# Volume(1h) label update for row 47 (ACDXExchange/ACXT now flagged "Worst in 24h").
$ws.Range("E47").Value = "46ACDXExchangeACXTWorstin24h"

$ws.Range("D48").Value = "'0.6797"
$ws.Range("D49").Value = "'0.1633"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("D51").Value = "'0.01010"
